# Skill.xlsx edit: add DamageCount / DamageIntervalTime columns (internal
# time for skill) between "CoolDownTime" (old col L) and the rest of the
# table; insert two new columns at L:M shifting everything from the old
# L onward two columns to the right, then populate the new columns and
# adjust the trailing "DefaultHitTime" column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two new columns (this also shifts cell formatting,
#        column widths and the data-validation sqref automatically). ---
$ws.Range("L1:M1").EntireColumn.Insert()

# Match the width Excel gives "DamageCount"/"DamageIntervalTime" (same as
# the neighbouring 20-wide columns I:K).
$refWidth = $ws.Range("I1").ColumnWidth
$ws.Range("L1:M1").ColumnWidth = $refWidth

# --- 2. Header / type / description rows for the new columns. ---
# Order matters: new shared-string entries are appended in first-use
# order, so write "DamageCount" / "伤害count" / "DamageIntervalTime" /
# "only be used..." in that exact sequence.
$ws.Range("L1").Value = "DamageCount"
$ws.Range("L10").Value = "伤害count"
$ws.Range("M1").Value = "DamageIntervalTime"
$ws.Range("M10").Value = "only be used when damage count > 1"

$ws.Range("L2").Value = "int"
$ws.Range("M2").Value = "float"

# --- 3. Data rows: every skill row gets DamageCount = 1 and
#        DamageIntervalTime = 0.3. ---
$ws.Range("L11:L46").Value = 1
$ws.Range("M11:M46").Value = 0.3

# --- 4. Adjust the trailing "DefaultHitTime" column (now column Y, was
#        W before the insert) for the new multi-hit timing. ---
# Rows 11-18: single default hit time 0.8 -> 0.6.
$ws.Range("Y11:Y18").Value = 0.6

# Row 19 stays 0.6 (unchanged).

# Rows 20,22,23: per-hit list "0.6,0.6,0.6" (first use of this string).
$ws.Range("Y20").Value = "0.6,0.6,0.6"
# Row 21: per-hit list "0.6,0.6,0.7" (second new string).
$ws.Range("Y21").Value = "0.6,0.6,0.7"
$ws.Range("Y22").Value = "0.6,0.6,0.6"
$ws.Range("Y23").Value = "0.6,0.6,0.6"

# Rows 24-28: same per-hit list, but these cells also get word-wrap
# turned on (new cell style).
$ws.Range("Y24").Value = "0.6,0.6,0.6"
$ws.Range("Y25").Value = "0.6,0.6,0.6"
$ws.Range("Y26").Value = "0.6,0.6,0.6"
$ws.Range("Y27").Value = "0.6,0.6,0.6"
$ws.Range("Y28").Value = "0.6,0.6,0.6"
$ws.Range("Y24:Y28").WrapText = $true

# Rows 29-46 stay 0.6 (unchanged).

# --- 5. View bookkeeping to match where the editor ended up. ---
$ws.Range("Y21").Select()
